$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 49, shifting existing rows 49:87 down to 50:88
$ws.Rows.Item(49).Insert()

# Populate new row 49 with the new weekly entry
$ws.Range("A49").Value = 10
$ws.Range("B49").Value = "Vega Modelo de Temuco"
$ws.Range("C49").Value = "La Araucanía"
$ws.Range("D49").Value = 45118
$ws.Range("E49").Value = 9
$ws.Range("F49").Value = 100112042
$ws.Range("G49").Value = "Locoto"
$ws.Range("H49").Value = "Sin especificar"
$ws.Range("I49").Value = "Primera"
$ws.Range("J49").Value = 180
$ws.Range("K49").Value = 3800
$ws.Range("L49").Value = 3800
$ws.Range("M49").Value = 3800
$ws.Range("N49").Value = '$/kilo'
$ws.Range("O49").Value = "Región de Arica y Parinacota"
$ws.Range("P49").Value = 3800
$ws.Range("Q49").Value = 1
$ws.Range("R49").Value = "Hortaliza"
